$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 190, shifting existing rows 190-254 down to 191-255.
$ws.Rows(190).Insert()

# Populate the newly inserted row 190 with the new market observation.
$ws.Range("A190").Value = 5
$ws.Range("B190").Value = "Macroferia Regional de Talca"
$ws.Range("C190").Value = "Maule"
$ws.Range("D190").Value = 44627
$ws.Range("E190").Value = 7
$ws.Range("F190").Value = 100112003
$ws.Range("G190").Value = "Ajo"
$ws.Range("H190").Value = "Chino"
$ws.Range("I190").Value = "Primera"
$ws.Range("J190").Value = 200
$ws.Range("K190").Value = 20000
$ws.Range("L190").Value = 20000
$ws.Range("M190").Value = 20000
$ws.Range("N190").Value = "$/malla 10 kilos"
$ws.Range("O190").Value = "China"
$ws.Range("P190").Value = 2000
$ws.Range("Q190").Value = 10
$ws.Range("R190").Value = "Hortaliza"
